$d = $word.ActiveDocument

[void]$d.Content.Find.Execute("My program files are ", $true, $false, $false, $false, $false, $true, 1, $false, "The application’s program files are ", 2)
[void]$d.Content.Find.Execute("re I prompt the user with a menu and read in the user’s input. I have a switch statement to determine what option the use", $true, $false, $false, $false, $false, $true, 1, $false, "re the user is prompted with a menu and the user’s input is read in. A switch statement determines what option the use", 2)
[void]$d.Content.Find.Execute(". I have", $true, $false, $false, $false, $false, $true, 1, $false, ". There is", 2)
[void]$d.Content.Find.Execute("rather handled the data as it was retrieved", $true, $false, $false, $false, $false, $true, 1, $false, "rather handle the data as it is retrieved", 2)
[void]$d.Content.Find.Execute("At the beginning of my program a Session object is constructed to connect to ", $true, $false, $false, $false, $false, $true, 1, $false, "At the beginning of the application a Session object is constructed with the user details to connect to ", 2)
[void]$d.Content.Find.Execute("ySQL. From there I query the ", $true, $false, $false, $false, $false, $true, 1, $false, "ySQL. From there the application queries the ", 2)
[void]$d.Content.Find.Execute(" ‘Experiment’. Once I switch to using the ", $true, $false, $false, $false, $false, $true, 1, $false, " ‘Experiment’. The application then switches to using the ", 2)
[void]$d.Content.Find.Execute("database I create ", $true, $false, $false, $false, $false, $true, 1, $false, "database and creates ", 2)
[void]$d.Content.Find.Execute(" into a separate function. One of the main things about my application is that when a user goes to insert an experiment or a run, it is all or nothing. They cannot break ", $true, $false, $false, $false, $false, $true, 1, $false, " into a separate function. One of the important rules within the application is that when a user goes to insert an experiment or a run, it is all or nothing. The user cannot break ", 2)
[void]$d.Content.Find.Execute("For inserting an experiment, I ask the user for the ", $true, $false, $false, $false, $false, $true, 1, $false, "For inserting an experiment, the application prompts the user for the ", 2)
[void]$d.Content.Find.Execute(", making sure there isn’t already an experiment in the table with that id. From there I ask the user for the rest of the meta data. There", $true, $false, $false, $false, $false, $true, 1, $false, ", making sure there isn’t already an experiment in the table with that id. There", 2)
[void]$d.Content.Find.Execute("so I ask for and enter the experiment meta data", $true, $false, $false, $false, $false, $true, 1, $false, "so the user is prompted for the experiment meta data", 2)
[void]$d.Content.Find.Execute(" and insert it all into the database", $true, $false, $false, $false, $false, $true, 1, $false, " and it all gets inserted into the database", 2)
[void]$d.Content.Find.Execute(" you must enter all the ", $true, $false, $false, $false, $false, $true, 1, $false, " the user must enter all the ", 2)
[void]$d.Content.Find.Execute("When inserting either a run parameter or run result, I use regex to make sure the value the user enters, matches with the type of the parameter. To look up information about the experiment, I ask the user for the ", $true, $false, $false, $false, $false, $true, 1, $false, "When inserting either a run parameter or run result, the application uses regex to make sure the value the user enters, matches with the type of the parameter. To look up information about the experiment, the application prompts the user for the ", 2)
[void]$d.Content.Find.Execute("and making sure that it exists in the database, I query for the meta data and print it, then I query the ", $true, $false, $false, $false, $false, $true, 1, $false, "and, making sure that it exists in the database, queries for the meta data to print and then queries the ", 2)
[void]$d.Content.Find.Execute(" table matching to the experimentId to print out ", $true, $false, $false, $false, $false, $true, 1, $false, " table, matching to the experimentId, to print out ", 2)
[void]$d.Content.Find.Execute("first I print out all possible runs that belong to the ", $true, $false, $false, $false, $false, $true, 1, $false, "first the application prompts the user by printing out all possible runs that belong to the ", 2)
[void]$d.Content.Find.Execute(" that the user provides. Then I let the user select which run", $true, $false, $false, $false, $false, $true, 1, $false, " that the user provides. Then the user selects which run", 2)
[void]$d.Content.Find.Execute(". For the experiment report it is nearly the same as fetching data about the experiment except I open an html file and print the data along with html tags for creating tables. The aggregate report option prompts the user to pick a list of ", $true, $false, $false, $false, $false, $true, 1, $false, ". For the experiment report it is nearly the same as fetching data about the experiment except the application creates an html file and prints the data along with html tags for creating tables. The aggregate report option prompts the user to pick from a list of ", 2)
[void]$d.Content.Find.Execute(" int or float. Once they’ve selected the parameter I ask for a min and max date and then print the ", $true, $false, $false, $false, $false, $true, 1, $false, " int or float type. Once they’ve selected the parameter the application prompts for a min and max date and then prints the ", 2)
[void]$d.Content.Find.Execute("I make sure the user doesn’t enter the same date for the start and end date. Lastly for the parameter search I ask the user for a parameterName and type and use a SQL query to retrieve the meta data of each experiment with that parameter, by joining", $true, $false, $false, $false, $false, $true, 1, $false, "The application makes sure the user doesn’t enter the same date for the start and end date. Lastly for the parameter search the application asks the user for a parameterName and type and uses a SQL query to retrieve the meta data of each experiment with that parameter, by joining", 2)
[void]$d.Content.Find.Execute("ParameterType table.", $true, $false, $false, $false, $false, $true, 1, $false, "ParameterType tables.", 2)

# Restore the _GoBack bookmark that originally sat inside "doesn't e|nter"
$bm = $d.Content
[void]$bm.Find.Execute("doesn’t e")
$bmRange = $d.Range($bm.End, $bm.End)
$d.Bookmarks.Add("_GoBack", $bmRange)
